# Update the "Forecast Comparison" sheet: insert a new "Week_Start_Date"
# column after "Week" and before "ASIN", shift the rest of the table over,
# shorten the week labels, tweak a handful of forecast numbers, and mark
# the holiday flag column as boolean values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column B ("Week_Start_Date"); everything from the old B
# onward (ASIN, MyForecast, Amazon Mean/P70/P80/P90 Forecast, Product
# Title, is_holiday_week) shifts one column to the right automatically.
$ws.Columns("B").Insert()

$ws.Range("B1").Value = "Week_Start_Date"

$weekStarts = @("2025-01-05","2025-01-12","2025-01-19","2025-01-26","2025-02-02","2025-02-09","2025-02-16","2025-02-23","2025-03-02","2025-03-09","2025-03-16","2025-03-23","2025-03-30","2025-04-06","2025-04-13","2025-04-20")
$weekLabels = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")

# D..H after the insert: MyForecast, Amazon Mean Forecast, Amazon P70
# Forecast, Amazon P80 Forecast, Amazon P90 Forecast.
$myForecast = @(6,8,8,8,8,8,8,9,8,8,8,8,8,8,7,7)
$meanForecast = @(6,6,6,7,7,7,7,7,7,7,7,7,7,6,6,6)
$p70Forecast = @(7,7,7,8,8,8,9,8,8,8,8,8,8,7,8,7)
$p80Forecast = @(9,10,9,11,11,11,12,12,11,11,11,12,11,10,11,10)
$p90Forecast = @(11,13,12,15,16,15,17,17,16,16,17,17,17,16,16,16)

for ($i = 0; $i -lt 16; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $weekLabels[$i]
    $ws.Range("B$r").Value = "'" + $weekStarts[$i]
    $ws.Cells.Item($r, 4).Value = $myForecast[$i]
    $ws.Cells.Item($r, 5).Value = $meanForecast[$i]
    $ws.Cells.Item($r, 6).Value = $p70Forecast[$i]
    $ws.Cells.Item($r, 7).Value = $p80Forecast[$i]
    $ws.Cells.Item($r, 8).Value = $p90Forecast[$i]
    $ws.Cells.Item($r, 10).Value = $false
}

# Update the "Summary" sheet metrics affected by the corrected forecast.
# These are stored as text (matching the rest of the column), so force the
# leading apostrophe to stop Excel from re-typing them as numbers/dates.
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B9").Value = "'126"
$summary.Range("B10").Value = "'64"
$summary.Range("B12").Value = "'9"
$summary.Range("B13").Value = "'2025-02-23"
$summary.Range("B14").Value = "'6"
$summary.Range("B15").Value = "'2025-01-05"
